# Regenerate s_val data to filter save games:
# Update columns B:E (and derived sum column G) for rows 2-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569)
    3 = @(0.1554434735375247, 0.3375848360084654, 0.1529057820181812, 0.4998867070740569)
    4 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569)
    5 = @(0.7287194209349384, 9.226618575922256, 3.082599426703578, 6.48142807727062)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]

    $sum = $vals[0] + $vals[1] + $vals[2] + $vals[3]
    $ws.Range("G$row").Value = $sum
}
